$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# New identifiers for this handback run
# ---------------------------------------------------------------------------
$oldGuid1 = "1db42d52-e708-4e55-bb66-94e619fe2ede"
$newGuid1 = "89e05de2-83ad-4e5e-a67d-1b6105ce6083"
$oldGuid2 = "5fdde428-dd8a-42e4-9056-d2d53ae6f2fb"
$newGuid2 = "ffffbd6e53cd-1566-4cd4-8811-bf6365f183dc"

$newFile1 = $newGuid1 + ".md"
$newFile2 = $newGuid2 + ".md"

$newZhCnXlf = $newGuid1 + ".12379d27f82a64c2b63665bf50a6c8754ebed751.zh-cn.xlf"
$newDeDeXlf = $newGuid1 + ".12379d27f82a64c2b63665bf50a6c8754ebed751.de-de.xlf"

$newLatestXliffDate   = "2016-08-23 13:04:23"
$newZhCnHandoffDate   = "2016-08-23 13:04:19"
$newZhCnHandbackDate  = "2016-08-23 13:04:36"
$newDeDeHandbackDate  = "2016-08-23 13:04:43"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("G2").Value = $newLatestXliffDate
$wsOverview.Range("G3").Value = $newLatestXliffDate

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/952db2e39b6d51b1a02a78fc7291ec356cac3c62/e2e/" + $oldGuid1 + ".md", "", "", "e2e\" + $newFile1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/952db2e39b6d51b1a02a78fc7291ec356cac3c62/e2e/" + $oldGuid2 + ".md", "", "", "e2e\" + $newFile2)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("I3").Value = $newFile2

$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("J2").Value = $newZhCnXlf
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("J3").Value = $newZhCnXlf

$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate
$wsZhCn.Range("K2").Value = $newZhCnHandbackDate
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/952db2e39b6d51b1a02a78fc7291ec356cac3c62/e2e/" + $oldGuid1 + ".md", "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a0dca30bc659feec65bda852e5a468580ca4e0f4/e2e/" + $oldGuid1 + ".md", "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/952db2e39b6d51b1a02a78fc7291ec356cac3c62/e2e/" + $oldGuid2 + ".md", "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a0dca30bc659feec65bda852e5a468580ca4e0f4/e2e/" + $oldGuid2 + ".md", "", "", $newFile2)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("I3").Value = $newFile2

$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("J2").Value = $newDeDeXlf
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("J3").Value = $newDeDeXlf

$wsDeDe.Range("H2").Value = $newLatestXliffDate
$wsDeDe.Range("H3").Value = $newLatestXliffDate

$wsDeDe.Range("K2").Value = $newDeDeHandbackDate
$wsDeDe.Range("K3").Value = $newDeDeHandbackDate

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/952db2e39b6d51b1a02a78fc7291ec356cac3c62/e2e/" + $oldGuid1 + ".md", "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4315494839918b7c1a46a3ea8922c254034698d7/e2e/" + $oldGuid1 + ".md", "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/952db2e39b6d51b1a02a78fc7291ec356cac3c62/e2e/" + $oldGuid2 + ".md", "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4315494839918b7c1a46a3ea8922c254034698d7/e2e/" + $oldGuid2 + ".md", "", "", $newFile2)

Write-Output "Handback report generated."
